$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" everywhere it
# appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the same
# underlying string).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    while ($found -ne $null) {
        $found.Value = "In Translation"
        $found = $used.FindNext($found)
    }
}

# Re-fit the Status-ish columns that held the now-shorter text so their
# stored width reflects the new content.
$wb.Worksheets.Item("Overview").Range("E:F").ColumnWidth = 13.4101845877511
$wb.Worksheets.Item("zh-cn").Range("C:C").ColumnWidth = 13.4101845877511
$wb.Worksheets.Item("de-de").Range("C:C").ColumnWidth = 13.4101845877511
